$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.113.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.176.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.01%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.69%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.173.90'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('E10').Value = '  -6.35%  '
$ws.Range('E11').Value = '  -5.78%  '
$ws.Range('E12').Value = '  -3.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000236'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.700.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.176.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.029.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('E19').Value = '  -4.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('E22').Value = '  -5.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.49%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.55%  '
$ws.Range('E30').Value = '  -6.79%  '
$ws.Range('E31').Value = '  -5.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.94%  '
$ws.Range('E33').Value = '  -4.33%  '
$ws.Range('E34').Value = '  -6.35%  '
$ws.Range('E35').Value = '  -6.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.79'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0706'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.62%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '403.67'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.09'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('E42').Value = '  -4.11%  '
$ws.Range('E43').Value = '  -6.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.813.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.82%  '
$ws.Range('E45').Value = '  -5.48%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.11'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.85%  '
$ws.Range('E51').Value = '  -2.21%  '
